$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the estimates for the "one supported retailer" automation row (row 5)
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0.3

# Match the % complete number format used by the other rows in column D
$ws.Range("D5").NumberFormat = "0%"

# Update the active selection to D5 (matches the selection recorded in the file)
$ws.Range("D5").Select()
